# Append the new gene-expression rows (GLS clustering export) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 3; A = 'Rv3723'; B = 3; C = 'lucA Rv3723'; D = 'FUNCTION: Required for the import of both fatty acids and cholesterol during growth in macrophages and in axenic culture. Facilitates the uptake of these lipids by stabilizing protein subunits of the Mce1 and Mce4 multi-subunit transporters, which transport fatty acids and cholesterol, respectively. Required for full virulence in vivo. {ECO:0000269|PubMed:28708968}.'; E = 38 }
    @{ Row = 4; A = 'Rv0174'; B = 2; C = 'mce1F Rv0174'; D = $null; E = 38 }
    @{ Row = 5; A = 'Rv0175'; B = 2; C = 'Rv0175'; D = $null; E = 38 }
    @{ Row = 6; A = 'Rv0170'; B = 2; C = 'mce1B Rv0170'; D = $null; E = 38 }
    @{ Row = 7; A = 'Rv2536'; B = 2; C = 'Rv2536'; D = $null; E = 38 }
    @{ Row = 8; A = 'Rv0172'; B = 2; C = 'mce1D Rv0172'; D = $null; E = 38 }
    @{ Row = 9; A = 'Rv0178'; B = 2; C = 'Rv0178'; D = $null; E = 38 }
    @{ Row = 10; A = 'Rv1016c'; B = 2; C = 'lpqT Rv1016c MTCY10G2.33'; D = $null; E = 38 }
    @{ Row = 11; A = 'Rv1405c'; B = 1; C = 'Rv1405c MTCY21B4.22c'; D = $null; E = 38 }
    @{ Row = 12; A = 'Rv0167'; B = 1; C = 'yrbE1A Rv0167'; D = $null; E = 38 }
    @{ Row = 13; A = 'Rv0513'; B = 1; C = 'Rv0513'; D = $null; E = 38 }
    @{ Row = 14; A = 'Rv0177'; B = 1; C = 'Rv0177'; D = $null; E = 38 }
    @{ Row = 15; A = 'Rv0168'; B = 1; C = 'yrbE1B Rv0168'; D = $null; E = 38 }
    @{ Row = 16; A = 'Rv0200'; B = 1; C = 'Rv0200'; D = $null; E = 38 }
    @{ Row = 17; A = 'Rv0173'; B = 1; C = 'lprK Rv0173'; D = $null; E = 38 }
    @{ Row = 18; A = 'Rv0171'; B = 1; C = 'mce1C Rv0171'; D = $null; E = 38 }
    @{ Row = 19; A = 'Rv0176'; B = 1; C = 'Rv0176'; D = $null; E = 38 }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    if ($null -ne $r.D) {
        $ws.Cells.Item($r.Row, 4).Value = $r.D
    } else {
        $ws.Cells.Item($r.Row, 4).Value = ""
    }
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}

Write-Output "Added $($newRows.Count) rows (now A1:E19)"
